# Realestate Update resale numbers 2024-01-13 14:09
# Append a new data row (row 54) to the CityResaleNum sheet with the
# 2024-01-13 14:09:04 resale-number snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 54

# Columns A ("Date") and D ("Week") hold values that look numeric/date-like
# ("2024-01-13", "01"), so force Text format before assigning the literal
# strings to avoid Excel's automatic date/number inference, then clear the
# formatting back to the sheet's default (unstyled) cell look - matching how
# every other data row in this sheet is stored.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2024-01-13"
$ws.Range("A$row").ClearFormats()

$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "01"
$ws.Range("D$row").ClearFormats()

$ws.Range("B$row").Value = "14:09:04"
$ws.Range("C$row").Value = "Saturday"

$ws.Range("E$row").Value = 137934
$ws.Range("F$row").Value = 142867
$ws.Range("G$row").Value = 171794
$ws.Range("H$row").Value = 148362
$ws.Range("I$row").Value = -1
$ws.Range("J$row").Value = 119523
$ws.Range("K$row").Value = 224987
$ws.Range("L$row").Value = 253270
$ws.Range("M$row").Value = 184888
$ws.Range("N$row").Value = 110434
$ws.Range("O$row").Value = 40967
$ws.Range("P$row").Value = 30908
$ws.Range("Q$row").Value = 73112
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 42558
$ws.Range("T$row").Value = -1
